{"js": "// Remove the \"Version / Release Date\" table-like paragraph and the\n// paragraph holding the version number / release date values from the\n// README (these two paragraphs sat between the title block and the\n// \"Exclusion of Liability\" heading).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  const isVersionHeaderRow =\n    text.indexOf(\"Version\") !== -1 && text.indexOf(\"Release Date\") !== -1;\n  const isVersionDateRow = /^\\s*202\\d\\s*\\/\\s*\\d{1,2}\\b/.test(text);\n  if (isVersionHeaderRow || isVersionDateRow) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Version / Release Date\" label paragraph and the paragraph\n# holding the actual version number / release date values from the\n# README (these two paragraphs sat between the title block and the\n# \"Exclusion of Liability\" heading).\n\n$doc = $word.ActiveDocument\n\n$keepGoing = $true\nwhile ($keepGoing) {\n    $keepGoing = $false\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        $isVersionHeaderRow = ($t -like \"*Version*\") -and ($t -like \"*Release Date*\")\n        $isVersionDateRow = $t -match \"^\\s*202\\d\\s*/\\s*\\d{1,2}\\b\"\n        if ($isVersionHeaderRow -or $isVersionDateRow) {\n            $p.Range.Delete()\n            $keepGoing = $true\n            break\n        }\n    }\n}\n"}
